# New Submission Synced: 2026-02-07 19:37:39
#
# The "JSS 3A" sheet holds one Google-Form-style submission in row 2.
# A new submission came in, so:
#   - row 2's "Admission No" (C2) is normalized from a text "23" to the
#     numeric value 23
#   - a new row 3 is appended with the new submission's data, where the
#     "Admission No" (C3) stays textual ("14") while the "AI Score" (D3)
#     is numeric (10)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3A")

# Normalize C2 ("Admission No") from inline text "23" to a real number.
$ws.Cells.Item(2, 3).Value = 23

# Append the new submission as row 3.
$ws.Cells.Item(3, 1).Value = "2026-02-07 19:37:39"
$ws.Cells.Item(3, 2).Value = "Muhammad dahiru idrisa "

# Admission No "14" must remain text (not auto-converted to a number).
# Force text with a leading apostrophe, then clear the resulting
# quote-prefix formatting so the cell keeps the workbook's default style.
$ws.Cells.Item(3, 3).Value = "'14"
$ws.Cells.Item(3, 3).ClearFormats()

$ws.Cells.Item(3, 4).Value = 10
